$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11.95318627383047
$ws.Range("D2").Value = 4.738517207037305
$ws.Range("E2").Value = 11.43271364712148
$ws.Range("F2").Value = 76.67680597230219
$ws.Range("G2").Value = 3.890156783053762
$ws.Range("J2").Value = 12.24735252397662
$ws.Range("K2").Value = 28.79401591188076
$ws.Range("L2").Value = 8.223468179656093
$ws.Range("M2").Value = 26.34677980467688

$ws.Range("C3").Value = 11.96268488544085
$ws.Range("D3").Value = 4.658036342586691
$ws.Range("E3").Value = 11.46466957590377
$ws.Range("F3").Value = 76.13546744708496
$ws.Range("G3").Value = 3.895297061369014
$ws.Range("J3").Value = 12.25686859821378
$ws.Range("K3").Value = 28.77485417353097
$ws.Range("L3").Value = 8.23207485600294
$ws.Range("M3").Value = 26.39698210732318

$ws.Range("C4").Value = 11.97099046658522
$ws.Range("D4").Value = 4.607588710844154
$ws.Range("E4").Value = 11.48555058742409
$ws.Range("F4").Value = 75.8107266946855
$ws.Range("G4").Value = 3.898611801271855
$ws.Range("J4").Value = 12.26363514847687
$ws.Range("K4").Value = 28.77231748875564
$ws.Range("L4").Value = 8.237668917598667
$ws.Range("M4").Value = 26.43441577521697

$ws.Range("C5").Value = 11.9749962203706
$ws.Range("D5").Value = 4.586781791182443
$ws.Range("E5").Value = 11.49437753053826
$ws.Range("F5").Value = 75.68037310963075
$ws.Range("G5").Value = 3.900002642265823
$ws.Range("J5").Value = 12.26662462221721
$ws.Range("K5").Value = 28.77360402789567
$ws.Range("L5").Value = 8.240026631702008
$ws.Range("M5").Value = 26.45132887293439

$ws.Range("C6").Value = 11.97569886227735
$ws.Range("D6").Value = 4.583312082489142
$ws.Range("E6").Value = 11.49586245735993
$ws.Range("F6").Value = 75.65884928292218
$ws.Range("G6").Value = 3.900236014983041
$ws.Range("J6").Value = 12.2671350310626
$ws.Range("K6").Value = 28.77395774538535
$ws.Range("L6").Value = 8.240422852197971
$ws.Range("M6").Value = 26.45423734652417

$ws.Range("C7").Value = 11.97104197569765
$ws.Range("D7").Value = 4.607309096102745
$ws.Range("E7").Value = 11.48566834280193
$ws.Range("F7").Value = 75.80896060503517
$ws.Range("G7").Value = 3.898630396209618
$ws.Range("J7").Value = 12.26367452627736
$ws.Range("K7").Value = 28.77232544670024
$ws.Range("L7").Value = 8.237700398035706
$ws.Range("M7").Value = 26.43463716038748

$ws.Range("C8").Value = 11.95594762078488
$ws.Range("D8").Value = 4.710987768098863
$ws.Range("E8").Value = 11.44347111441135
$ws.Range("F8").Value = 76.48860255429575
$ws.Range("G8").Value = 3.891896344270079
$ws.Range("J8").Value = 12.25044184025257
$ws.Range("K8").Value = 28.78549388039293
$ws.Range("L8").Value = 8.226371672953553
$ws.Range("M8").Value = 26.36271590107244

$ws.Range("C9").Value = 11.94600160600826
$ws.Range("D9").Value = 4.90567241814093
$ws.Range("E9").Value = 11.37067747693396
$ws.Range("F9").Value = 77.87928418740348
$ws.Range("G9").Value = 3.879941027849453
$ws.Range("J9").Value = 12.23183115041962
$ws.Range("K9").Value = 28.88447471154484
$ws.Range("L9").Value = 8.206600496006871
$ws.Range("M9").Value = 26.2742552979998

$ws.Range("C10").Value = 11.95070510815055
$ws.Range("D10").Value = 5.042885895874075
$ws.Range("E10").Value = 11.32320701003661
$ws.Range("F10").Value = 78.93245667796492
$ws.Range("G10").Value = 3.87190809844647
$ws.Range("J10").Value = 12.2226460867359
$ws.Range("K10").Value = 29.00157997439504
$ws.Range("L10").Value = 8.193548791223812
$ws.Range("M10").Value = 26.24147715102843

$ws.Range("C11").Value = 11.95545480597762
$ws.Range("D11").Value = 5.103932905953672
$ws.Range("E11").Value = 11.30290449167163
$ws.Range("F11").Value = 79.41751798541635
$ws.Range("G11").Value = 3.868414239095245
$ws.Range("J11").Value = 12.21944517536811
$ws.Range("K11").Value = 29.06440398518939
$ws.Range("L11").Value = 8.187927937543856
$ws.Range("M11").Value = 26.23358533601553

$ws.Range("C12").Value = 11.9576284282675
$ws.Range("D12").Value = 5.12684370186903
$ws.Range("E12").Value = 11.29540130598714
$ws.Range("F12").Value = 79.60197764590642
$ws.Range("G12").Value = 3.867114073089269
$ws.Range("J12").Value = 12.21837384369933
$ws.Range("K12").Value = 29.08955665716785
$ws.Range("L12").Value = 8.185844707161804
$ws.Range("M12").Value = 26.23160735250853

$ws.Range("C13").Value = 11.95714362831803
$ws.Range("D13").Value = 5.121918770460399
$ws.Range("E13").Value = 11.29700903967239
$ws.Range("F13").Value = 79.56221768280707
$ws.Range("G13").Value = 3.867393072204786
$ws.Range("J13").Value = 12.21859830912798
$ws.Range("K13").Value = 29.08407917446996
$ws.Range("L13").Value = 8.186291358695531
$ws.Range("M13").Value = 26.23198839312058

$ws.Range("C14").Value = 11.95562611792724
$ws.Range("D14").Value = 5.105821977430892
$ws.Range("E14").Value = 11.30228349852364
$ws.Range("F14").Value = 79.43267847954628
$ws.Range("G14").Value = 3.868306816067002
$ws.Range("J14").Value = 12.21935421344415
$ws.Range("K14").Value = 29.06644607264569
$ws.Range("L14").Value = 8.187755643106327
$ws.Range("M14").Value = 26.23340235012038

$ws.Range("C15").Value = 11.95474542388313
$ws.Range("D15").Value = 5.095935088564847
$ws.Range("E15").Value = 11.30553831590886
$ws.Range("F15").Value = 79.35343080133018
$ws.Range("G15").Value = 3.868869485314282
$ws.Range("J15").Value = 12.21983556765073
$ws.Range("K15").Value = 29.05582237614768
$ws.Range("L15").Value = 8.188658447388825
$ws.Range("M15").Value = 26.23440005777728

$ws.Range("C16").Value = 11.95044720419118
$ws.Range("D16").Value = 5.038867933619851
$ws.Range("E16").Value = 11.32455975308498
$ws.Range("F16").Value = 78.90087001285286
$ws.Range("G16").Value = 3.872139642112385
$ws.Range("J16").Value = 12.22287496204538
$ws.Range("K16").Value = 28.99766559674566
$ws.Range("L16").Value = 8.193922474306095
$ws.Range("M16").Value = 26.24213425165858

$ws.Range("C17").Value = 11.94847873777487
$ws.Range("D17").Value = 5.003501051205453
$ws.Range("E17").Value = 11.33655909855357
$ws.Range("F17").Value = 78.62471032117237
$ws.Range("G17").Value = 3.874186723697996
$ws.Range("J17").Value = 12.22499004056158
$ws.Range("K17").Value = 28.96442765007908
$ws.Range("L17").Value = 8.197232660617175
$ws.Range("M17").Value = 26.24867767276498

$ws.Range("C18").Value = 11.94759223159697
$ws.Range("D18").Value = 4.983030322140642
$ws.Range("E18").Value = 11.34358246786046
$ws.Range("F18").Value = 78.46643725463144
$ws.Range("G18").Value = 3.875379255364614
$ws.Range("J18").Value = 12.22629856385853
$ws.Range("K18").Value = 28.94620995246562
$ws.Range("L18").Value = 8.1991663904566
$ws.Range("M18").Value = 26.25310191536663

$ws.Range("C19").Value = 11.94733428119432
$ws.Range("D19").Value = 4.97607747632918
$ws.Range("E19").Value = 11.34598138197639
$ws.Range("F19").Value = 78.41294847046059
$ws.Range("G19").Value = 3.875785625691362
$ws.Range("J19").Value = 12.22675739884812
$ws.Range("K19").Value = 28.94019661288569
$ws.Range("L19").Value = 8.199826243289428
$ws.Range("M19").Value = 26.25471330997582

$ws.Range("C20").Value = 11.94866285648173
$ws.Range("D20").Value = 5.007279296552839
$ws.Range("E20").Value = 11.33526916229588
$ws.Range("F20").Value = 78.65404992096362
$ws.Range("G20").Value = 3.873967246394113
$ws.Range("J20").Value = 12.2247553649282
$ws.Range("K20").Value = 28.96787282453043
$ws.Range("L20").Value = 8.196877203358149
$ws.Range("M20").Value = 26.24791273080733

$ws.Range("C21").Value = 11.95606167332926
$ws.Range("D21").Value = 5.110555666885474
$ws.Range("E21").Value = 11.30072924973662
$ws.Range("F21").Value = 79.47070683513108
$ws.Range("G21").Value = 3.868037807725771
$ws.Range("J21").Value = 12.21912836337607
$ws.Range("K21").Value = 29.07158845962257
$ws.Range("L21").Value = 8.187324320689424
$ws.Range("M21").Value = 26.23295960628787

$ws.Range("C22").Value = 11.96308279343577
$ws.Range("D22").Value = 5.176845919865202
$ws.Range("E22").Value = 11.27923298524357
$ws.Range("F22").Value = 80.00893244710441
$ws.Range("G22").Value = 3.86429587245722
$ws.Range("J22").Value = 12.21627149155526
$ws.Range("K22").Value = 29.14730829664315
$ws.Range("L22").Value = 8.181344689004089
$ws.Range("M22").Value = 26.22907684338238

$ws.Range("C23").Value = 11.95913567203583
$ws.Range("D23").Value = 5.141578864851929
$ws.Range("E23").Value = 11.29060762958361
$ws.Range("F23").Value = 79.72128728325136
$ws.Range("G23").Value = 3.86628087567022
$ws.Range("J23").Value = 12.21772108914903
$ws.Range("K23").Value = 29.10617322975905
$ws.Range("L23").Value = 8.184512077835599
$ws.Range("M23").Value = 26.23060997953032

$ws.Range("C24").Value = 11.94857885267808
$ws.Range("D24").Value = 5.005571580469357
$ws.Range("E24").Value = 11.3358519537335
$ws.Range("F24").Value = 78.64078394055008
$ws.Range("G24").Value = 3.874066423364117
$ws.Range("J24").Value = 12.22486117356909
$ws.Range("K24").Value = 28.96631248433372
$ws.Range("L24").Value = 8.197037810068993
$ws.Range("M24").Value = 26.24825649810924

$ws.Range("C25").Value = 11.9465842809705
$ws.Range("D25").Value = 4.853988199574939
$ws.Range("E25").Value = 11.3893103387989
$ws.Range("F25").Value = 77.49725514481592
$ws.Range("G25").Value = 3.883042608074474
$ws.Range("J25").Value = 12.23607866932814
$ws.Range("K25").Value = 28.84987986193373
$ws.Range("L25").Value = 8.211689091635746
$ws.Range("M25").Value = 26.29253891553178

